$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 77001576
$ws.Range("I33").Value = 111113340
$ws.Range("J33").Value = 250094.5
$ws.Range("K33").Value = 111113340
$ws.Range("L33").Value = 250094.5
$ws.Range("M33").Value = -111113111
$ws.Range("N33").Value = -250552.5

$ws.Range("H40").Value = 1410.2727
$ws.Range("I40").Value = 1307.875
$ws.Range("J40").Value = 1683.3334
$ws.Range("K40").Value = 1307.875
$ws.Range("L40").Value = 1683.3334
$ws.Range("M40").Value = -1132.875
$ws.Range("N40").Value = -2033.3334

$ws.Range("H113").Value = 238900.23
$ws.Range("I113").Value = 311008.62
$ws.Range("J113").Value = 2909.0908
$ws.Range("K113").Value = 311008.62
$ws.Range("L113").Value = 2909.0908
$ws.Range("M113").Value = -307754.62
$ws.Range("N113").Value = -9417.0908

$ws.Range("H138").Value = 3420.9822
$ws.Range("I138").Value = 1816.3889
$ws.Range("J138").Value = 6309.25
$ws.Range("K138").Value = 5449.1667
$ws.Range("L138").Value = 18927.75
$ws.Range("M138").Value = -309.1666999999998
$ws.Range("N138").Value = -29207.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2326.75
$ws.Range("I2").Value = 2618.6
$ws.Range("J2").Value = 1451.2
$ws.Range("K2").Value = 2618.6
$ws.Range("L2").Value = 1451.2
$ws.Range("M2").Value = -2505.6
$ws.Range("N2").Value = -1677.2

$ws.Range("H45").Value = 1463.8182
$ws.Range("I45").Value = 1309.4445
$ws.Range("K45").Value = 1309.4445
$ws.Range("M45").Value = -932.4445000000001

$ws.Range("H61").Value = 2195.3333
$ws.Range("I61").Value = 1128.5714
$ws.Range("J61").Value = 2634.5881
$ws.Range("K61").Value = 1128.5714
$ws.Range("L61").Value = 2634.5881
$ws.Range("M61").Value = -916.5714
$ws.Range("N61").Value = -3058.5881

$ws.Range("H63").Value = 2231.6328
$ws.Range("I63").Value = 2217.554
$ws.Range("J63").Value = 2440
$ws.Range("K63").Value = 2217.554
$ws.Range("L63").Value = 2440
$ws.Range("M63").Value = -1531.554
$ws.Range("N63").Value = -3812

$ws.Range("H66").Value = 2231.6328
$ws.Range("I66").Value = 2217.554
$ws.Range("J66").Value = 2440
$ws.Range("K66").Value = 11087.77
$ws.Range("L66").Value = 12200
$ws.Range("M66").Value = -7655.77
$ws.Range("N66").Value = -19064

$ws.Range("H97").Value = 1480
$ws.Range("I97").Value = 1436.0769
$ws.Range("J97").Value = 1670.3334
$ws.Range("K97").Value = 1436.0769
$ws.Range("L97").Value = 1670.3334
$ws.Range("M97").Value = -940.0769
$ws.Range("N97").Value = -2662.3334

$ws.Range("H102").Value = 1755.875
$ws.Range("I102").Value = 1692.6897
$ws.Range("J102").Value = 2366.6667
$ws.Range("K102").Value = 1692.6897
$ws.Range("L102").Value = 2366.6667
$ws.Range("M102").Value = -70.6896999999999
$ws.Range("N102").Value = -5610.6667

$ws.Range("H116").Value = 2326.75
$ws.Range("I116").Value = 2618.6
$ws.Range("J116").Value = 1451.2
$ws.Range("K116").Value = 2618.6
$ws.Range("L116").Value = 1451.2
$ws.Range("M116").Value = -324.5999999999999
$ws.Range("N116").Value = -6039.2

$ws.Range("H136").Value = 2195.3333
$ws.Range("I136").Value = 1128.5714
$ws.Range("J136").Value = 2634.5881
$ws.Range("K136").Value = 3385.7142
$ws.Range("L136").Value = 7903.7643
$ws.Range("M136").Value = -835.7142000000003
$ws.Range("N136").Value = -13003.7643

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2326.75
$ws.Range("I3").Value = 2618.6
$ws.Range("J3").Value = 1451.2
$ws.Range("K3").Value = 2618.6
$ws.Range("L3").Value = 1451.2
$ws.Range("M3").Value = -2504.6
$ws.Range("N3").Value = -1679.2

$ws.Range("H22").Value = 443.57144
$ws.Range("I22").Value = 328.18182
$ws.Range("J22").Value = 866.6667
$ws.Range("K22").Value = 328.18182
$ws.Range("L22").Value = 866.6667
$ws.Range("M22").Value = -155.18182

$ws.Range("H99").Value = 1888.7368
$ws.Range("I99").Value = 1959.0667
$ws.Range("J99").Value = 1625
$ws.Range("K99").Value = 1959.0667
$ws.Range("L99").Value = 1625
$ws.Range("M99").Value = -461.0667000000001
$ws.Range("N99").Value = -4621

$ws.Range("H105").Value = 1264426.2
$ws.Range("I105").Value = 2842147.8
$ws.Range("J105").Value = 2249.1
$ws.Range("K105").Value = 2842147.8
$ws.Range("L105").Value = 2249.1
$ws.Range("M105").Value = -2840400.8
$ws.Range("N105").Value = -5743.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3298.8948
$ws.Range("I16").Value = 2781.6667
$ws.Range("J16").Value = 3537.6155
$ws.Range("K16").Value = 2781.6667
$ws.Range("L16").Value = 3537.6155
$ws.Range("M16").Value = -2494.6667
$ws.Range("N16").Value = -4111.6155

$ws.Range("H22").Value = 27778280
$ws.Range("I22").Value = 35714696
$ws.Range("J22").Value = 824.75
$ws.Range("K22").Value = 35714696
$ws.Range("L22").Value = 824.75
$ws.Range("M22").Value = -35714346
$ws.Range("N22").Value = -1524.75

$ws.Range("H31").Value = 7814012.5
$ws.Range("I31").Value = 31251056
$ws.Range("K31").Value = 31251056
$ws.Range("M31").Value = -31250761

$ws.Range("H34").Value = 7814012.5
$ws.Range("I34").Value = 31251056
$ws.Range("K34").Value = 31251056
$ws.Range("M34").Value = -31250854

$ws.Range("H86").Value = 325612.25
$ws.Range("I86").Value = 689074.9
$ws.Range("J86").Value = 2534.3333
$ws.Range("K86").Value = 689074.9
$ws.Range("L86").Value = 2534.3333
$ws.Range("M86").Value = -687951.9
$ws.Range("N86").Value = -4780.3333

$ws.Range("H89").Value = 325612.25
$ws.Range("I89").Value = 689074.9
$ws.Range("J89").Value = 2534.3333
$ws.Range("K89").Value = 3445374.5
$ws.Range("L89").Value = 12671.6665
$ws.Range("M89").Value = -3439758.5
$ws.Range("N89").Value = -23903.6665

$ws.Range("H105").Value = 1269.5264
$ws.Range("I105").Value = 955.4545000000001
$ws.Range("J105").Value = 1701.375
$ws.Range("K105").Value = 955.4545000000001
$ws.Range("L105").Value = 1701.375
$ws.Range("M105").Value = 791.5454999999999
$ws.Range("N105").Value = -5195.375

$ws.Range("H107").Value = 935.625
$ws.Range("I107").Value = 369.90475
$ws.Range("J107").Value = 2015.6364
$ws.Range("K107").Value = 369.90475
$ws.Range("L107").Value = 2015.6364
$ws.Range("M107").Value = 1550.09525
$ws.Range("N107").Value = -5855.6364

$ws.Range("H113").Value = 3298.8948
$ws.Range("I113").Value = 2781.6667
$ws.Range("J113").Value = 3537.6155
$ws.Range("K113").Value = 2781.6667
$ws.Range("L113").Value = 3537.6155
$ws.Range("M113").Value = -611.6667000000002
$ws.Range("N113").Value = -7877.6155

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1242.84
$ws.Range("I68").Value = 886.85364
$ws.Range("J68").Value = 1672.1177
$ws.Range("K68").Value = 2660.56092
$ws.Range("L68").Value = 5016.3531
$ws.Range("M68").Value = -1849.56092
$ws.Range("N68").Value = -6638.3531

$ws.Range("H71").Value = 1242.84
$ws.Range("I71").Value = 886.85364
$ws.Range("J71").Value = 1672.1177
$ws.Range("K71").Value = 7981.682760000001
$ws.Range("L71").Value = 15049.0593
$ws.Range("M71").Value = -3925.682760000001
$ws.Range("N71").Value = -23161.0593

$ws.Range("H107").Value = 427.34375
$ws.Range("J107").Value = 1258.5834
$ws.Range("L107").Value = 3775.7502
$ws.Range("N107").Value = -7615.7502

$ws.Range("H113").Value = 3247162.8
$ws.Range("I113").Value = 427.46667
$ws.Range("J113").Value = 10204453
$ws.Range("K113").Value = 1282.40001
$ws.Range("L113").Value = 30613359
$ws.Range("M113").Value = 887.5999899999999
$ws.Range("N113").Value = -30617699

$ws.Range("H131").Value = 10639172
$ws.Range("I131").Value = 35714750
$ws.Range("J131").Value = 1046.3939
$ws.Range("K131").Value = 107144250
$ws.Range("L131").Value = 3139.1817
$ws.Range("M131").Value = -107139210
$ws.Range("N131").Value = -13219.1817

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3727.6052
$ws.Range("J93").Value = 1766.1666
$ws.Range("L93").Value = 1766.1666
$ws.Range("N93").Value = -4262.1666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1349.9
$ws.Range("I107").Value = 2000
$ws.Range("J107").Value = 1187.375
$ws.Range("K107").Value = 6000
$ws.Range("L107").Value = 3562.125
$ws.Range("M107").Value = -4080
$ws.Range("N107").Value = -7402.125

Write-Host "Edit complete"